# Remove the space in player names throughout the match-history sheet,
# clear the "ここまで初期データ" marker cell, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("対戦履歴")

$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in 2, 3) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -eq $null) { continue }

        if ($val -like "吉谷*") {
            $cell.Value = "吉谷悠"
        } elseif ($val -like "石井山*") {
            $cell.Value = "石井山拓登"
        } elseif ($val -like "森岡*") {
            $cell.Value = "森岡凜太郎"
        }
    }
}

# Clear the "ここまで初期データ" marker in D33 but keep its style.
$ws.Cells.Item(33, 4).ClearContents()

# Move the view back to the top and select D12 (matches the saved state).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D12").Select()
